$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 2790.3635
$ws.Range("J62").Value = 3014.2856
$ws.Range("L62").Value = 3014.2856
$ws.Range("N62").Value = -4262.2856

# Row 65
$ws.Range("H65").Value = 2790.3635
$ws.Range("J65").Value = 3014.2856
$ws.Range("L65").Value = 15071.428
$ws.Range("N65").Value = -21311.428

# Row 69
$ws.Range("H69").Value = 0
$ws.Range("I69").Value = 0
$ws.Range("K69").Value = 0
$ws.Range("M69").ClearContents()

# Row 72
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = 0
$ws.Range("K72").Value = 0
$ws.Range("M72").ClearContents()

# Row 80
$ws.Range("H80").Value = 825
$ws.Range("I80").Value = 987.5
$ws.Range("J80").Value = 500
$ws.Range("K80").Value = 2962.5
$ws.Range("L80").Value = 1500
$ws.Range("M80").Value = -1964.5
$ws.Range("N80").Value = -3496

# Row 83
$ws.Range("H83").Value = 825
$ws.Range("I83").Value = 987.5
$ws.Range("J83").Value = 500
$ws.Range("K83").Value = 8887.5
$ws.Range("L83").Value = 4500
$ws.Range("M83").Value = -3895.5
$ws.Range("N83").Value = -14484

# Row 86
$ws.Range("H86").Value = 10283.1
$ws.Range("I86").Value = 5136.2
$ws.Range("J86").Value = 15430
$ws.Range("K86").Value = 5136.2
$ws.Range("L86").Value = 15430
$ws.Range("M86").Value = -4013.2
$ws.Range("N86").Value = -17676

# Row 89
$ws.Range("H89").Value = 10283.1
$ws.Range("I89").Value = 5136.2
$ws.Range("J89").Value = 15430
$ws.Range("K89").Value = 25681
$ws.Range("L89").Value = 77150
$ws.Range("M89").Value = -20065
$ws.Range("N89").Value = -88382

# Row 104
$ws.Range("H104").Value = 264.66666
$ws.Range("I104").Value = 264.66666
$ws.Range("K104").Value = 793.9999799999999
$ws.Range("M104").Value = 953.0000200000001

# Row 112
$ws.Range("H112").Value = 2633693.5
$ws.Range("I112").Value = 10001060
$ws.Range("J112").Value = 2491.3572
$ws.Range("K112").Value = 30003180
$ws.Range("L112").Value = 7474.071599999999
$ws.Range("M112").Value = -30002072
$ws.Range("N112").Value = -9690.071599999999

# Row 127
$ws.Range("H127").Value = 1754.4546
$ws.Range("I127").Value = 1229.9
$ws.Range("K127").Value = 3689.7
$ws.Range("M127").Value = 1270.3

# Row 137
$ws.Range("H137").Value = 1176.8
$ws.Range("I137").Value = 1043.2307
$ws.Range("J137").Value = 2045
$ws.Range("K137").Value = 3129.6921
$ws.Range("L137").Value = 6135
$ws.Range("M137").Value = -579.6921000000002
$ws.Range("N137").Value = -11235

# Row 138
$ws.Range("H138").Value = 2126
$ws.Range("I138").Value = 848.6667
$ws.Range("J138").Value = 5000
$ws.Range("K138").Value = 2546.0001
$ws.Range("L138").Value = 15000
$ws.Range("M138").Value = 2593.9999
$ws.Range("N138").Value = -25280

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 4195.45
$ws.Range("I32").Value = 4009.353
$ws.Range("K32").Value = 4009.353
$ws.Range("M32").Value = -3722.353

# Row 61
$ws.Range("H61").Value = 3515
$ws.Range("I61").Value = 3515
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3515
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3303
$ws.Range("N61").ClearContents()

# Row 106
$ws.Range("H106").Value = 31333
$ws.Range("J106").Value = 31333
$ws.Range("L106").Value = 31333
$ws.Range("N106").Value = -33857

# Row 109
$ws.Range("H109").Value = 30374.5
$ws.Range("J109").Value = 30374.5
$ws.Range("L109").Value = 30374.5
$ws.Range("N109").Value = -33148.5

# Row 132
$ws.Range("H132").Value = 1188.9231
$ws.Range("I132").Value = 1095.8
$ws.Range("K132").Value = 3287.4
$ws.Range("M132").Value = -757.3999999999996

# Row 136
$ws.Range("H136").Value = 3515
$ws.Range("I136").Value = 3515
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 10545
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -7995
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 134
$ws.Range("H134").Value = 6900.4443
$ws.Range("I134").Value = 6448.3477
$ws.Range("J134").Value = 9500
$ws.Range("K134").Value = 19345.0431
$ws.Range("L134").Value = 28500
$ws.Range("M134").Value = -16810.0431
$ws.Range("N134").Value = -33570

# Row 135
$ws.Range("H135").Value = 49999
$ws.Range("J135").Value = 49999
$ws.Range("L135").Value = 49999
$ws.Range("N135").Value = -60139

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 3535.7576
$ws.Range("I31").Value = 2672.5
$ws.Range("K31").Value = 2672.5
$ws.Range("M31").Value = -2377.5

# Row 34
$ws.Range("H34").Value = 3535.7576
$ws.Range("I34").Value = 2672.5
$ws.Range("K34").Value = 2672.5
$ws.Range("M34").Value = -2470.5

# Row 45
$ws.Range("H45").Value = 0
$ws.Range("J45").Value = 0
$ws.Range("L45").Value = 0
$ws.Range("N45").ClearContents()

# Row 68
$ws.Range("H68").Value = 21237.867
$ws.Range("I68").Value = 4821
$ws.Range("K68").Value = 4821
$ws.Range("M68").Value = -4072

# Row 71
$ws.Range("H71").Value = 21237.867
$ws.Range("I71").Value = 4821
$ws.Range("K71").Value = 14463
$ws.Range("M71").Value = -10719

# Row 82
$ws.Range("H82").Value = 22000
$ws.Range("J82").Value = 22000
$ws.Range("L82").Value = 22000
$ws.Range("N82").Value = -22722

# Row 85
$ws.Range("H85").Value = 22000
$ws.Range("J85").Value = 22000
$ws.Range("L85").Value = 22000
$ws.Range("N85").Value = -24496

# Row 132
$ws.Range("H132").Value = 1225.2354
$ws.Range("I132").Value = 864.3125
$ws.Range("K132").Value = 2592.9375
$ws.Range("M132").Value = -62.9375

$ws = $wb.Worksheets.Item("CUL")
# Row 80
$ws.Range("H80").Value = 10949.25
$ws.Range("J80").Value = 13999.667
$ws.Range("L80").Value = 41999.001
$ws.Range("N80").Value = -43871.001

# Row 83
$ws.Range("H83").Value = 10949.25
$ws.Range("J83").Value = 13999.667
$ws.Range("L83").Value = 125997.003
$ws.Range("N83").Value = -135357.003

# Row 121
$ws.Range("H121").Value = 1012.6667
$ws.Range("J121").Value = 1086.1
$ws.Range("L121").Value = 3258.3
$ws.Range("N121").Value = -5878.299999999999

# Row 140
$ws.Range("H140").Value = 1817.8182
$ws.Range("J140").Value = 5500
$ws.Range("L140").Value = 16500
$ws.Range("N140").Value = -26860

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 5000
$ws.Range("I70").Value = 5000
$ws.Range("K70").Value = 5000
$ws.Range("M70").Value = -4730

# Row 73
$ws.Range("H73").Value = 5000
$ws.Range("I73").Value = 5000
$ws.Range("K73").Value = 5000
$ws.Range("M73").Value = -4064

# Row 132
$ws.Range("H132").Value = 1803.6
$ws.Range("I132").Value = 1892.8889
$ws.Range("K132").Value = 5678.6667
$ws.Range("M132").Value = -3148.6667

$ws = $wb.Worksheets.Item("LTW")
# Row 55
$ws.Range("H55").Value = 440
$ws.Range("I55").Value = 400
$ws.Range("J55").Value = 466.66666
$ws.Range("K55").Value = 400
$ws.Range("L55").Value = 466.66666
$ws.Range("M55").Value = -227
$ws.Range("N55").Value = -812.66666

# Row 64
$ws.Range("H64").Value = 44833.332
$ws.Range("J64").Value = 29750
$ws.Range("L64").Value = 29750
$ws.Range("N64").Value = -30200

# Row 67
$ws.Range("H67").Value = 44833.332
$ws.Range("J67").Value = 29750
$ws.Range("L67").Value = 29750
$ws.Range("N67").Value = -31310

# Row 108
$ws.Range("H108").Value = 58496.5
$ws.Range("J108").Value = 58496.5
$ws.Range("L108").Value = 58496.5
$ws.Range("N108").Value = -66176.5

# Row 132
$ws.Range("H132").Value = 6129.4375
$ws.Range("I132").Value = 4461
$ws.Range("K132").Value = 13383
$ws.Range("M132").Value = -10853

# Row 136
$ws.Range("H136").Value = 2530.6155
$ws.Range("I136").Value = 2415.0908
$ws.Range("J136").Value = 3166
$ws.Range("K136").Value = 7245.2724
$ws.Range("L136").Value = 9498
$ws.Range("M136").Value = -4695.2724
$ws.Range("N136").Value = -14598

$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 4838.5557
$ws.Range("I62").Value = 4966.1665
$ws.Range("J62").Value = 4583.3335
$ws.Range("K62").Value = 4966.1665
$ws.Range("L62").Value = 4583.3335
$ws.Range("M62").Value = -4342.1665
$ws.Range("N62").Value = -5831.3335

# Row 65
$ws.Range("H65").Value = 4838.5557
$ws.Range("I65").Value = 4966.1665
$ws.Range("J65").Value = 4583.3335
$ws.Range("K65").Value = 24830.8325
$ws.Range("L65").Value = 22916.6675
$ws.Range("M65").Value = -21710.8325
$ws.Range("N65").Value = -29156.6675

# Row 101
$ws.Range("H101").Value = 59646.332
$ws.Range("J101").Value = 59646.332
$ws.Range("L101").Value = 59646.332
$ws.Range("N101").Value = -66136.33199999999

# Row 104
$ws.Range("H104").Value = 23666.666
$ws.Range("J104").Value = 23666.666
$ws.Range("L104").Value = 23666.666
$ws.Range("N104").Value = -30654.666

# Row 126
$ws.Range("H126").Value = 1474.5
$ws.Range("I126").Value = 1474.5
$ws.Range("K126").Value = 4423.5
$ws.Range("M126").Value = -1953.5

# Row 132
$ws.Range("H132").Value = 1681.7333
$ws.Range("I132").Value = 1575.0834
$ws.Range("K132").Value = 4725.2502
$ws.Range("M132").Value = -2195.2502

Write-Host "Applied all updates"